$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "48.280.73"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.36%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.511.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.67%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.03%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "

# Row 7
$ws.Range("E7").Value = "  +1.09%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.542"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.10%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.59%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.02%  "

# Row 12
$ws.Range("E12").Value = "  +1.16%  "

# Row 13
$ws.Range("E13").Value = "  +0.10%  "

# Row 14
$ws.Range("E14").Value = "  -0.35%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.902.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.57%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.506.21"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.51%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.845"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "48.132.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.36%  "

# Row 19
$ws.Range("E19").Value = "  -1.56%  "

# Row 20
$ws.Range("E20").Value = "  +2.46%  "

# Row 21
$ws.Range("E21").Value = "  +0.38%  "

# Row 22
$ws.Range("E22").Value = "  +0.28%  "

# Row 23
$ws.Range("B23").Value = "BitcoinCash"
$ws.Range("C23").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "278.84"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +12.90%  "

# Row 24
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.42%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "

# Row 26
$ws.Range("E26").Value = "  -0.07%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.80"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.16%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.21"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.98%  "

# Row 29
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.03%  "

# Row 30
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.76%  "

# Row 31
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.35"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.93%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.78%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.74"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.18%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.35"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.45%  "

# Row 35
$ws.Range("B35").Value = "FirstDigitalUSD"
$ws.Range("C35").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.10%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0785"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.35%  "

# Row 37
$ws.Range("E37").Value = "  +0.03%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.68"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "

# Row 39
$ws.Range("E39").Value = "  -1.08%  "

# Row 40
$ws.Range("E40").Value = "  -0.09%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "122.35"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.75%  "

# Row 42
$ws.Range("E42").Value = "  +0.24%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.48"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.19%  "

# Row 44
$ws.Range("E44").Value = "  +1.93%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.021.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.28%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.20"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.74%  "

# Row 47
$ws.Range("E47").Value = "  +3.16%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.35%  "

# Row 49
$ws.Range("E49").Value = "  -0.01%  "

# Row 50
$ws.Range("E50").Value = "  -0.64%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.78"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.06%  "
